# LC changes/fixed for better alignment
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated B, D, E, F values for rows 2-15
$data = @{
    2  = @{ B = 1267.090909090909;  D = 1267.216666666666; E = -0.1257575757574614; F = -0.1257575757574614 }
    3  = @{ B = 1374.781609195402;  D = 1374.607142857143; E = 0.1744663382596627;  F = 0.3002239140171241 }
    4  = @{ B = 1408.801652892562;  D = 1408.576923076923; E = 0.2247298156391935;  F = 0.05026347737953074 }
    5  = @{ B = 1642.663366336634;  D = 1643.236111111111; E = -0.572744774477087;  F = -0.7974745901162805 }
    6  = @{ B = 1680.396226415095;  D = 1680.652777777778; E = -0.2565513626832399; F = 0.3161934117938472 }
    7  = @{ B = 1710.015625;        D = 1709.916666666666; E = 0.09895833333371229; F = 0.3555096960169521 }
    8  = @{ B = 1775.159090909091;  D = 1774.535714285714; E = 0.623376623376771;   F = 0.5244182900430587 }
    9  = @{ B = 2088.276923076923;  D = 2087.994047619048; E = 0.2828754578754342;  F = -0.3405011655013368 }
    10 = @{ B = 2113.0875;          D = 2112.735714285714; E = 0.3517857142855974;  F = 0.06891025641016313 }
    11 = @{ B = 2373.765957446808;  D = 2373.626984126984; E = 0.1389733198243448;  F = -0.2128123944612526 }
    12 = @{ B = 2519.5;             D = 2519.253246753247; E = 0.246753246753542;   F = 0.1077799269291972 }
    13 = @{ B = 2554.012195121951;  D = 2553.422077922078; E = 0.5901171998734753;  F = 0.3433639531199333 }
    14 = @{ B = 2669.344262295082;  D = 2669.448051948052; E = -0.1037896529701356; F = -0.693906852843611 }
    15 = @{ B = 3119;               D = 3118.714285714286; E = 0.2857142857142208;  F = 0.3895039386843564 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals.B
    $ws.Cells.Item($row, 4).Value = $vals.D
    $ws.Cells.Item($row, 5).Value = $vals.E
    $ws.Cells.Item($row, 6).Value = $vals.F
}

# Remove rows 16, 17, 18 (old tail rows no longer present)
$ws.Range("A16:F18").EntireRow.Delete() | Out-Null
